$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$targetCells = @('D2','E2','D3','D4','E4','D5','E5','D6','E6','D7','E7','D8','E8','D9','E9','D10','E10','D11','E11','D12','E12','D13','E13','D14','E14','D15','E15','D16','E16','D17','E17','E18','D19','E19','D20','E20','E21','E22','D23','E23','D24','E24','D25','E25','E26','D39','E39','D40','E40','D41','E41','D42','E42','D43','E43','D44','E44','D45','E45','D46','E46','E47','D48','E48','D50','E50','E51')
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '307.59'
$ws.Range('E2').Value = '-2.52%'
$ws.Range('D3').Value = '37.87'
$ws.Range('D4').Value = '5.045'
$ws.Range('E4').Value = '-1.90%'
$ws.Range('D5').Value = '0.07892'
$ws.Range('E5').Value = '-3.52%'
$ws.Range('D6').Value = '1.990'
$ws.Range('E6').Value = '1.07%'
$ws.Range('D7').Value = '4.366'
$ws.Range('E7').Value = '2.97%'
$ws.Range('D8').Value = '8.213'
$ws.Range('E8').Value = '-0.07%'
$ws.Range('D9').Value = '3.177'
$ws.Range('E9').Value = '1.25%'
$ws.Range('D10').Value = '0.9252'
$ws.Range('E10').Value = '-0.10%'
$ws.Range('D11').Value = '0.1272'
$ws.Range('E11').Value = '-9.33%'
$ws.Range('D12').Value = '0.1900'
$ws.Range('E12').Value = '-4.02%'
$ws.Range('D13').Value = '0.08705'
$ws.Range('E13').Value = '-3.72%'
$ws.Range('D14').Value = '0.03460'
$ws.Range('E14').Value = '-1.44%'
$ws.Range('D15').Value = '0.09727'
$ws.Range('E15').Value = '-1.02%'
$ws.Range('D16').Value = '0.001392'
$ws.Range('E16').Value = '0.09%'
$ws.Range('D17').Value = '0.005943'
$ws.Range('E17').Value = '-0.49%'
$ws.Range('E18').Value = '-3.06%'
$ws.Range('D19').Value = '0.3434'
$ws.Range('E19').Value = '-0.92%'
$ws.Range('D20').Value = '0.1299'
$ws.Range('E20').Value = '-0.32%'
$ws.Range('E21').Value = '5.33%'
$ws.Range('E22').Value = '3.58%'
$ws.Range('D23').Value = '0.04327'
$ws.Range('E23').Value = '-1.15%'
$ws.Range('D24').Value = '0.001224'
$ws.Range('E24').Value = '0.16%'
$ws.Range('D25').Value = '0.004595'
$ws.Range('E25').Value = '-4.03%'
$ws.Range('E26').Value = '176.66%'
$ws.Range('D39').Value = '0.02254'
$ws.Range('E39').Value = '3.32%'
$ws.Range('D40').Value = '0.04994'
$ws.Range('E40').Value = '-3.68%'
$ws.Range('D41').Value = '0.007579'
$ws.Range('E41').Value = '0.36%'
$ws.Range('D42').Value = '0.009872'
$ws.Range('E42').Value = '0.84%'
$ws.Range('D43').Value = '0.1356'
$ws.Range('E43').Value = '-1.36%'
$ws.Range('D44').Value = '0.002094'
$ws.Range('E44').Value = '-1.60%'
$ws.Range('D45').Value = '0.008539'
$ws.Range('E45').Value = '-6.45%'
$ws.Range('D46').Value = '0.00006416'
$ws.Range('E46').Value = '0.44%'
$ws.Range('E47').Value = '0.20%'
$ws.Range('D48').Value = '0.003002'
$ws.Range('E48').Value = '8.66%'
$ws.Range('D50').Value = '0.00002103'
$ws.Range('E50').Value = '0.20%'
$ws.Range('E51').Value = '0.20%'
